$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the "Price" (D) cells being updated so that
# numeric-looking strings (e.g. "1.00", "7.00", "0.504") are stored
# verbatim as text, exactly like the rest of the price column,
# instead of being auto-coerced into trimmed numeric values by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.835.96"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "2.945.14"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "593.03"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "147.30"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "2.942.56"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").Value = "0.504"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  +5.46%  "
$ws.Range("D12").Value = "0.438"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("D14").Value = "32.54"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "3.435.39"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "62.838.54"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").Value = "6.67"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "2.946.22"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "438.68"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "13.42"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").Value = "0.664"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").Value = "7.00"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "11.18"
$ws.Range("E24").Value = "  +3.51%  "
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").Value = "2.12"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").Value = "11.81"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("E30").Value = "  +6.00%  "
$ws.Range("D31").Value = "2.60"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").Value = "0.0000101"
$ws.Range("E32").Value = "  +15.89%  "
$ws.Range("D33").Value = "26.32"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").Value = "5.61"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "3.03"
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("D39").Value = "49.67"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").Value = "2.02"
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.118"
$ws.Range("E41").Value = "  -3.48%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "8.44"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").Value = "0.278"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "39.10"
$ws.Range("E44").Value = "  -7.02%  "
$ws.Range("D45").Value = "2.705.94"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "135.21"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").Value = "357.34"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "22.67"
$ws.Range("E51").Value = "  -3.52%  "
